$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet2 ("NewAccount"): refresh the emailed-credentials list.
#    Remove the stale hyperlink on A2, rewrite its value, then append two new
#    rows (A3, A4) with their own hyperlinks, and finally re-link A2.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Range("A2").Value = "lakshmi@mailinator.com"

$ws2.Range("A3").Value = "Jasu@mailinator.com"
$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:Jasu@mailinator.com")
$ws2.Range("A3").Style = "Hyperlink"

$ws2.Range("A4").Value = "Sathish@mailinator.com"
$ws2.Hyperlinks.Add($ws2.Range("A4"), "mailto:Sathish@mailinator.com")
$ws2.Range("A4").Style = "Hyperlink"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:lakshmi@mailinator.com")
$ws2.Range("A2").Style = "Hyperlink"

$ws2.Range("B2").Select()

# ---------------------------------------------------------------------------
# 2. New sheet "AccountCreationData" with the account-creation test data.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "AccountCreationData"

# ---- header row -------------------------------------------------------
$headerRange = $ws4.Range("A1:J1")
$ws4.Range("A1").Value = "Email"
$ws4.Range("B1").Value = "Gender"
$ws4.Range("C1").Value = "FirstName"
$ws4.Range("D1").Value = "LastName"
$ws4.Range("E1").Value = "SetPassword"
$ws4.Range("F1").Value = "Day"
$ws4.Range("G1").Value = "Month"
$ws4.Range("H1").Value = "Year"
$ws4.Range("I1").Value = "Newsletter"
$ws4.Range("J1").Value = "OptinOffer"

$fullHeaderRange = $ws4.Range("A1:E1,I1:O1")
$fullHeaderRange.Font.Bold = $true
$fullHeaderRange.Borders.LineStyle = 1

$dateHeaderRange = $ws4.Range("F1:H1")
$dateHeaderRange.Font.Bold = $true
$dateHeaderRange.Borders.LineStyle = 1
$dateHeaderRange.NumberFormat = "@"

# ---- data rows ----------------------------------------------------------
# Row 2 - Lakshmi
$ws4.Range("B2").Value = "Female"
$ws4.Range("C2").Value = "Lakshmi"
$ws4.Range("D2").Value = "S"
$ws4.Range("I2").Value = "Yes"
$ws4.Range("J2").Value = "No"

$ws4.Range("E2").NumberFormat = "General"
$ws4.Range("E2").Value = "Lakshmi@7"
$ws4.Hyperlinks.Add($ws4.Range("E2"), "mailto:Lakshmi@7")
$ws4.Range("E2").Style = "Hyperlink"

# Row 3 - Jaswanthika
$ws4.Range("A3").Value = "Jasu6@mailinator.com"
$ws4.Hyperlinks.Add($ws4.Range("A3"), "mailto:Jasu6@mailinator.com")
$ws4.Range("A3").Style = "Hyperlink"
$ws4.Range("B3").Value = "female"
$ws4.Range("C3").Value = "Jaswanthika"
$ws4.Range("D3").Value = "S"
$ws4.Range("I3").Value = "No"
$ws4.Range("J3").Value = "Yes"

$ws4.Range("E3").NumberFormat = "General"
$ws4.Range("E3").Value = "Jasu@31"
$ws4.Hyperlinks.Add($ws4.Range("E3"), "mailto:Jasu@31")
$ws4.Range("E3").Style = "Hyperlink"

# Row 4 - Sathish
$ws4.Range("A4").Value = "Sathish6@mailinator.com"
$ws4.Hyperlinks.Add($ws4.Range("A4"), "mailto:Sathish6@mailinator.com")
$ws4.Range("A4").Style = "Hyperlink"
$ws4.Range("B4").Value = "male"
$ws4.Range("C4").Value = "Sathish"
$ws4.Range("D4").Value = "R"
$ws4.Range("I4").Value = "no"
$ws4.Range("J4").Value = "no"

$ws4.Range("E4").NumberFormat = "General"
$ws4.Range("E4").Value = "sathish@14"
$ws4.Hyperlinks.Add($ws4.Range("E4"), "mailto:sathish@14")
$ws4.Range("E4").Style = "Hyperlink"

$ws4.Range("A2").Value = "lakshmi6@mailinator.com"
$ws4.Hyperlinks.Add($ws4.Range("A2"), "mailto:lakshmi6@mailinator.com")
$ws4.Range("A2").Style = "Hyperlink"

# ---- Day / Month / Year columns are stored as text -----------------------
$dateDataRange = $ws4.Range("F2:H4")
$dateDataRange.NumberFormat = "@"

$ws4.Range("F2").Value = "7"
$ws4.Range("G2").Value = "9"
$ws4.Range("H2").Value = "1990"

$ws4.Range("F3").Value = "31"
$ws4.Range("G3").Value = "5"
$ws4.Range("H3").Value = "2021"

$ws4.Range("F4").Value = "14"
$ws4.Range("G4").Value = "6"
$ws4.Range("H4").Value = "1987"

# ---- layout ---------------------------------------------------------------
$ws4.Columns.Item(5).ColumnWidth = 12.28515625
$ws4.PageSetup.Orientation = 1

$ws4.Activate()
$ws4.Range("A4").Select()

Write-Output "done"
